$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Fecha de Factura o Documento" column (D2:D10): the RPA flow now
# writes a human readable date string instead of the raw date-serial text.
$ws.Range("D2:D10").Value = "28/10/19"

# Apply the built-in Hyperlink look (underline, themed link colour, text
# number format) to the e-mail column (I2:I10) without leaving an actual
# clickable hyperlink object behind - add then immediately remove the
# hyperlinks, which leaves the formatting/style in place.
$emailRange = $ws.Range("I2:I10")
foreach ($cell in $emailRange.Cells) {
    $addr = $cell.Value2
    $ws.Hyperlinks.Add($cell, ("mailto:" + $addr), "", "", $addr) | Out-Null
}
$ws.Hyperlinks.Delete()

# Widen column C slightly (30.29 -> ~31.57 chars) to match the new layout.
$ws.Columns.Item(3).ColumnWidth = 30.6

# Move the active selection to J14, matching the final cursor position.
$ws.Range("J14").Select() | Out-Null
